# Apply updated cryptocurrency price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.007.63'
$ws.Range('E2').Value = '  -0.12%  '
$ws.Range('D3').Value = '2.304.80'
$ws.Range('E3').Value = '  +0.00%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''306.25'
$ws.Range('E5').Value = '  +1.80%  '
$ws.Range('D6').Value = '''98.15'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').Value = '''0.513'
$ws.Range('E7').Value = '  -1.55%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = '''0.508'
$ws.Range('E9').Value = '  -1.91%  '
$ws.Range('D10').Value = '''35.93'
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('E11').Value = '  +0.11%  '
$ws.Range('D12').Value = '''18.29'
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('D13').Value = '''0.118'
$ws.Range('E13').Value = '  +1.10%  '
$ws.Range('D14').Value = '''6.81'
$ws.Range('E14').Value = '  -0.89%  '
$ws.Range('D15').Value = '2.663.77'
$ws.Range('E15').Value = '  +0.03%  '
$ws.Range('D16').Value = '2.300.08'
$ws.Range('E16').Value = '  -0.59%  '
$ws.Range('D17').Value = '''0.786'
$ws.Range('E17').Value = '  +0.00%  '
$ws.Range('D18').Value = '42.941.46'
$ws.Range('E18').Value = '  -0.07%  '
$ws.Range('E19').Value = '  -5.73%  '
$ws.Range('D20').Value = '0.0₃0905'
$ws.Range('E20').Value = '  -0.38%  '
$ws.Range('D21').Value = '''6.05'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('D22').Value = '''68.02'
$ws.Range('E22').Value = '  -0.29%  '
$ws.Range('D23').Value = '''236.71'
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('E24').Value = '  -1.45%  '
$ws.Range('D25').Value = '''2.51'
$ws.Range('E25').Value = '  +3.51%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('E27').Value = '  -0.37%  '
$ws.Range('D28').Value = '''25.55'
$ws.Range('E28').Value = '  +3.39%  '
$ws.Range('D29').Value = '''165.93'
$ws.Range('E29').Value = '  -1.59%  '
$ws.Range('E30').Value = '  +0.69%  '
$ws.Range('B31').Value = 'Cosmos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D31').Value = '''9.08'
$ws.Range('E31').Value = '  -0.68%  '
$ws.Range('B32').Value = 'InjectiveProtocol'
$ws.Range('C32').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D32').Value = '''33.40'
$ws.Range('E32').Value = '  +1.51%  '
$ws.Range('E33').Value = '  +0.07%  '
$ws.Range('D35').Value = '''5.04'
$ws.Range('E35').Value = '  -2.73%  '
$ws.Range('D36').Value = '''17.12'
$ws.Range('E36').Value = '  -5.33%  '
$ws.Range('E37').Value = '  -1.23%  '
$ws.Range('D38').Value = '''0.0692'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('E39').Value = '  -1.35%  '
$ws.Range('E40').Value = '  -1.03%  '
$ws.Range('E41').Value = '  -1.20%  '
$ws.Range('D42').Value = '''2.75'
$ws.Range('E42').Value = '  -0.17%  '
$ws.Range('D43').Value = '2.009.55'
$ws.Range('E43').Value = '  -0.12%  '
$ws.Range('E44').Value = '  -1.96%  '
$ws.Range('E45').Value = '  -0.90%  '
$ws.Range('D46').Value = '''17.81'
$ws.Range('E46').Value = '  +2.03%  '
$ws.Range('E47').Value = '  -3.52%  '
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('E49').Value = '  +3.68%  '
$ws.Range('D50').Value = '''53.76'
$ws.Range('E50').Value = '  -1.33%  '
$ws.Range('D51').Value = '2.531.35'
$ws.Range('E51').Value = '  +0.10%  '
